# Weekly update: insert a new week's worth of records (2 rows) at the top of the
# existing data block (rows 471-472), pushing all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 471; existing rows 471..524 shift to 473..526.
$ws.Rows.Item(471).Resize(2).Insert()

# Populate the first new row (471) - "Primera" quality record for the new week.
$ws.Cells.Item(471, 1).Value  = 9
$ws.Cells.Item(471, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(471, 3).Value  = "Metropolitana"
$ws.Cells.Item(471, 4).Value  = 45154
$ws.Cells.Item(471, 5).Value  = 13
$ws.Cells.Item(471, 6).Value  = 100112017
$ws.Cells.Item(471, 7).Value  = "Apio"
$ws.Cells.Item(471, 8).Value  = "Americana (o)"
$ws.Cells.Item(471, 9).Value  = "Primera"
$ws.Cells.Item(471, 10).Value = 70
$ws.Cells.Item(471, 11).Value = 7000
$ws.Cells.Item(471, 12).Value = 8000
$ws.Cells.Item(471, 13).Value = 7500
$ws.Cells.Item(471, 14).Value = "`$/docena de matas"
$ws.Cells.Item(471, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(471, 16).Value = 1250
$ws.Cells.Item(471, 17).Value = 6
$ws.Cells.Item(471, 18).Value = "Hortaliza"

# Populate the second new row (472) - "Segunda" quality record for the new week.
$ws.Cells.Item(472, 1).Value  = 9
$ws.Cells.Item(472, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(472, 3).Value  = "Metropolitana"
$ws.Cells.Item(472, 4).Value  = 45154
$ws.Cells.Item(472, 5).Value  = 13
$ws.Cells.Item(472, 6).Value  = 100112017
$ws.Cells.Item(472, 7).Value  = "Apio"
$ws.Cells.Item(472, 8).Value  = "Americana (o)"
$ws.Cells.Item(472, 9).Value  = "Segunda"
$ws.Cells.Item(472, 10).Value = 52
$ws.Cells.Item(472, 11).Value = 6000
$ws.Cells.Item(472, 12).Value = 6000
$ws.Cells.Item(472, 13).Value = 6000
$ws.Cells.Item(472, 14).Value = "`$/docena de matas"
$ws.Cells.Item(472, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(472, 16).Value = 1000
$ws.Cells.Item(472, 17).Value = 6
$ws.Cells.Item(472, 18).Value = "Hortaliza"

# Ensure the date cells keep the same numeric date format used elsewhere in column D.
$ws.Cells.Item(471, 4).NumberFormat = $ws.Cells.Item(473, 4).NumberFormat
$ws.Cells.Item(472, 4).NumberFormat = $ws.Cells.Item(473, 4).NumberFormat
